# Auto-generated Word COM-interop script
# Applies the commit diff: font TimesNewToman -> Times New Roman,
# and rewrites the 'Celestial Choreography' essay into the 'History' essay.

$d = $word.ActiveDocument

# --- Replace paragraphs whose text content changed (font baked into the new XML) ---
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="44"/>
        </w:rPr>
        <w:t>History: A Tapestry of Human Endeavors</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(1).Range.InsertXML($xml)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t>Dr</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve"> Henry Longfellow</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(2).Range.InsertXML($xml)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>hlongfellow@schoolhistory</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>org</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(3).Range.InsertXML($xml)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>History has always enticed us, capturing our imagination with its tales of triumphs, tragedies, remarkable achievements, and lessons learned</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Within its vast expanse, we unearth the foundations of civilizations, the intricacies of diverse cultures, and the evolution of human thought and action</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> It serves as a mosaic composed of countless individual stories, each contributing its unique hue to the broader canvas of our shared heritage</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> As we delve into the depths of historical studies, we assume the role of explorers, unraveling the intricate threads that bind us to past events and shedding light on the origins of our present world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>The tapestry of history is woven with the endeavors of countless individuals, both renowned and forgotten</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Kings and queens, generals and statesmen, revolutionaries and reformers, artists and thinkers - each has left their indelible mark on the canvas of human progress</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> But history would be incomplete if it solely focused on the exploits of the few; it is the collective actions and contributions of ordinary individuals, the unsung heroes of history, that truly shape its course</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> As students of history, we have the privilege of delving into the lives of these individuals, understanding their motivations, struggles, and triumphs, and appreciating the profound impact they have had on our world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Furthermore, history offers us a lens through which we can examine ourselves, our values, and the challenges we confront in our ever-changing society</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> By understanding the past, we gain perspective on the present and can make more informed decisions about the future</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> History compels us to confront uncomfortable truths, to grapple with difficult lessons, and to cultivate empathy for those who have come before us</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> In this way, history serves as a guide, helping us navigate the complexities of human existence and equipping us with the knowledge and wisdom to shape a better future for ourselves and generations to come</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(5).Range.InsertXML($xml)

# --- Update font on paragraphs whose text is unchanged ('Summary' heading) ---
$d.Paragraphs(6).Range.Font.Name = "Times New Roman"

# --- Replace the final 'Summary' paragraph AND append a new trailing empty paragraph
#     in one shot (InsertXML over paragraph 7's range, injecting two <w:p> elements) ---
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>History, in its vastness, is a narrative of human experiences, triumphs and challenges, that weaves together the fabric of our past to illuminate the present and guide the future</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> Through history, we gain insights into the lives of individuals who have shaped our world, the lessons learned from their endeavors, and the complex dynamics that have influenced the course of civilization</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> History is not merely a collection of facts; it is a tapestry of stories that holds the keys to our </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>collective identity, values, and aspirations</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> By studying history, we embark on a transformative journey that deepens our understanding of ourselves and our place within the grand tapestry of human existence</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(7).Range.InsertXML($xml)

# --- Safety net: make sure every remaining run uses Times New Roman ---
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
}
